# Agiledoc.xlsx edit: update "Original Estimate (hours)" column on the
# Sprint Backlog sheet, set a completion-date note on the last task row,
# and leave the Sprint Backlog sheet/tab active & selected (matching the
# author's last interaction before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# Make this the active sheet (drives workbookView.activeTab + tabSelected)
$ws.Activate()

# Revise the "Original Estimate (hours)" values in column I
$ws.Range("I2").Value  = 6
$ws.Range("I3").Value  = 8
$ws.Range("I4").Value  = 12
$ws.Range("I5").Value  = 10
$ws.Range("I6").Value  = 8
$ws.Range("I7").Value  = 12
$ws.Range("I8").Value  = 6
$ws.Range("I9").Value  = 8
$ws.Range("I11").Value = 12

# Add the completion-date note for the last task row
$ws.Range("E12").Value = "Sept 28th"

# Leave the selection where the author left it
$ws.Range("G19").Select()
